# Add a new task row ("another test task") with its time value, extending
# the task/time table by one row (B12:C12), matching the style/number
# format already used by the other task rows in columns B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "another test task"
$ws.Range("C12").Value = 0.003611111111111111

# Match the "0.00" number format used by the other time cells in column C
# (reuses the existing style rather than creating a new one).
$ws.Range("C12").NumberFormat = "0.00"
